$d = $word.ActiveDocument

# Replace the ${descripción} placeholder (kept as a literal placeholder text)
$d.Content.Find.Execute("`${descripción}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "`${descripción}", 2)

# Replace the ${jefe_mant} placeholder with the actual name
$d.Content.Find.Execute("`${jefe_mant}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ING. PAULINO LÓPEZ MODESTO", 2)

# Replace the ${jefe_dep} placeholder with the actual name (note trailing space)
$d.Content.Find.Execute("`${jefe_dep}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ING. SILVANO MARTÍNEZ HERNÁNDEZ ", 2)
